# Update cryptos list values to reflect latest market data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextCell $ws 'D2' '69.039.92'
Set-TextCell $ws 'E2' '  +1.56%  '
Set-TextCell $ws 'D3' '3.759.37'
Set-TextCell $ws 'E3' '  -0.58%  '
Set-TextCell $ws 'E4' '  -0.36%  '
Set-TextCell $ws 'D5' '629.22'
Set-TextCell $ws 'E5' '  +3.58%  '
Set-TextCell $ws 'D6' '165.38'
Set-TextCell $ws 'E6' '  +1.41%  '
Set-TextCell $ws 'D7' '3.758.45'
Set-TextCell $ws 'E7' '  -0.53%  '
Set-TextCell $ws 'E8' '  +0.08%  '
Set-TextCell $ws 'D9' '0.519'
Set-TextCell $ws 'E9' '  +0.77%  '
Set-TextCell $ws 'D10' '0.159'
Set-TextCell $ws 'E10' '  +0.69%  '
Set-TextCell $ws 'D11' '0.458'
Set-TextCell $ws 'E11' '  +2.40%  '
Set-TextCell $ws 'D12' '6.80'
Set-TextCell $ws 'E12' '  -0.04%  '
Set-TextCell $ws 'E13' '  -0.45%  '
Set-TextCell $ws 'D14' '34.88'
Set-TextCell $ws 'E14' '  +0.03%  '
Set-TextCell $ws 'D15' '4.391.60'
Set-TextCell $ws 'E15' '  -0.38%  '
Set-TextCell $ws 'D16' '3.724.93'
Set-TextCell $ws 'E16' '  -2.10%  '
Set-TextCell $ws 'D17' '68.961.60'
Set-TextCell $ws 'E17' '  +1.58%  '
Set-TextCell $ws 'D18' '17.58'
Set-TextCell $ws 'E18' '  -2.81%  '
Set-TextCell $ws 'D19' '0.114'
Set-TextCell $ws 'E19' '  -1.18%  '
Set-TextCell $ws 'D20' '6.98'
Set-TextCell $ws 'E20' '  -0.43%  '
Set-TextCell $ws 'D21' '464.32'
Set-TextCell $ws 'E21' '  +0.86%  '
Set-TextCell $ws 'D22' '9.54'
Set-TextCell $ws 'E22' '  +0.85%  '
Set-TextCell $ws 'D23' '0.703'
Set-TextCell $ws 'E23' '  +1.93%  '
Set-TextCell $ws 'D24' '82.91'
Set-TextCell $ws 'E24' '  -0.23%  '
Set-TextCell $ws 'E25' '  -1.59%  '
Set-TextCell $ws 'B26' 'Fetch.AI'
Set-TextCell $ws 'C26' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell $ws 'D26' '2.16'
Set-TextCell $ws 'E26' '  +4.13%  '
Set-TextCell $ws 'B27' 'InternetComputer(DFINITY)'
Set-TextCell $ws 'C27' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell $ws 'D27' '11.94'
Set-TextCell $ws 'E27' '  +0.47%  '
Set-TextCell $ws 'D28' '10.07'
Set-TextCell $ws 'E28' '  +1.58%  '
Set-TextCell $ws 'E29' '  -0.01%  '
Set-TextCell $ws 'D30' '3.907.21'
Set-TextCell $ws 'E30' '  -0.30%  '
Set-TextCell $ws 'B31' 'ImmutableX'
Set-TextCell $ws 'C31' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell $ws 'D31' '2.25'
Set-TextCell $ws 'E31' '  +4.08%  '
Set-TextCell $ws 'B32' 'PancakeSwap'
Set-TextCell $ws 'C32' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell $ws 'D32' '2.66'
Set-TextCell $ws 'E32' '  +2.34%  '
Set-TextCell $ws 'D33' '7.10'
Set-TextCell $ws 'E33' '  -1.06%  '
Set-TextCell $ws 'D34' '28.43'
Set-TextCell $ws 'E34' '  -1.92%  '
Set-TextCell $ws 'D35' '0.172'
Set-TextCell $ws 'E35' '  +16.20%  '
Set-TextCell $ws 'E36' '  +0.13%  '
Set-TextCell $ws 'D37' '3.709.15'
Set-TextCell $ws 'D38' '8.93'
Set-TextCell $ws 'E38' '  +0.14%  '
Set-TextCell $ws 'E39' '  +1.45%  '
Set-TextCell $ws 'D40' '3.30'
Set-TextCell $ws 'E40' '  +3.62%  '
Set-TextCell $ws 'E41' '  -0.71%  '
Set-TextCell $ws 'D42' '0.962'
Set-TextCell $ws 'E42' '  -1.49%  '
Set-TextCell $ws 'D43' '0.999'
Set-TextCell $ws 'E43' '  -0.09%  '
Set-TextCell $ws 'E44' '  -0.15%  '
Set-TextCell $ws 'D45' '156.12'
Set-TextCell $ws 'E45' '  +2.37%  '
Set-TextCell $ws 'D46' '43.36'
Set-TextCell $ws 'E46' '  +0.65%  '
Set-TextCell $ws 'B47' 'ONDO'
Set-TextCell $ws 'C47' 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextCell $ws 'D47' '1.41'
Set-TextCell $ws 'E47' '  +1.23%  '
Set-TextCell $ws 'B48' 'Stacks'
Set-TextCell $ws 'C48' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell $ws 'D48' '1.93'
Set-TextCell $ws 'E48' '  +5.39%  '
Set-TextCell $ws 'B49' 'OKB'
Set-TextCell $ws 'C49' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell $ws 'D49' '46.80'
Set-TextCell $ws 'E49' '  -0.67%  '
Set-TextCell $ws 'B50' 'TheGraph'
Set-TextCell $ws 'C50' 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextCell $ws 'D50' '0.295'
Set-TextCell $ws 'E50' '  +0.89%  '
Set-TextCell $ws 'B51' 'Cosmos'
Set-TextCell $ws 'C51' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell $ws 'D51' '8.32'
Set-TextCell $ws 'E51' '  +0.14%  '
